# Apply weekly refresh of Fruta/Hortaliza data:
# The values in columns D, J, K, L, M, P (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) for rows 2-31 get
# redistributed across rows according to the mapping below (after_row -> source_row).
# All other columns (A,B,C,E,F,G,H,I,N,O,Q,R) are identical across rows and remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: key = destination row, value = source row (values copied from source row's
# original content into destination row).
$rowMap = @{
    2  = 8
    3  = 2
    4  = 9
    5  = 14
    6  = 11
    7  = 17
    8  = 31
    9  = 4
    10 = 12
    11 = 15
    12 = 24
    13 = 19
    14 = 3
    15 = 27
    16 = 26
    17 = 6
    18 = 29
    19 = 28
    20 = 21
    21 = 7
    22 = 22
    23 = 16
    24 = 18
    25 = 20
    26 = 25
    27 = 10
    28 = 13
    29 = 23
    30 = 30
    31 = 5
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the original values for the affected columns before any writes happen,
# since rows reference each other as sources.
$original = @{}
foreach ($r in 2..31) {
    foreach ($c in $cols) {
        $original["$c$r"] = $ws.Range("$c$r").Value2
    }
}

foreach ($destRow in ($rowMap.Keys | Sort-Object)) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $original["$c$srcRow"]
    }
}
